$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.451284646987915
$ws.Range("B1").Value = 1.50369393825531
$ws.Range("C1").Value = 1.608202695846558
$ws.Range("D1").Value = 2.262841939926147
$ws.Range("E1").Value = 3.802537441253662
